# Applies the scheduled market-data refresh described in the commit
# "chore: update Sheets via scheduled runner".
#
# Each Leve row stores scraped Universalis market-board data in columns
# H:N (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ). This run refreshes those columns with newly scraped
# prices for the rows the scraper touched this cycle. A couple of rows
# (ALC!48 and ALC!56) lost HQ market data entirely, so their HQ profit
# cell (N) is cleared and the NQ profit cell (M) is updated instead.

$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 2023.5
$ws.Range("J39").Value = 1594.2
$ws.Range("L39").Value = 4782.6
$ws.Range("N39").Value = -5374.6
$ws.Range("H48").Value = 4332.6665
$ws.Range("I48").Value = 4332.6665
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 12997.9995
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -12705.9995
$ws.Range("N48").ClearContents()
$ws.Range("H56").Value = 4332.6665
$ws.Range("I56").Value = 4332.6665
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 12997.9995
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -12463.9995
$ws.Range("N56").ClearContents()
$ws.Range("H137").Value = 2349.84
$ws.Range("I137").Value = 1609.5333
$ws.Range("J137").Value = 3460.3
$ws.Range("K137").Value = 4828.5999
$ws.Range("L137").Value = 10380.9
$ws.Range("M137").Value = -2278.5999
$ws.Range("N137").Value = -15480.9
$ws.Range("H138").Value = 3199.0566
$ws.Range("J138").Value = 3526.122
$ws.Range("L138").Value = 10578.366
$ws.Range("N138").Value = -20858.366

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1287.56
$ws.Range("I2").Value = 1241.7894
$ws.Range("J2").Value = 1432.5
$ws.Range("K2").Value = 1241.7894
$ws.Range("L2").Value = 1432.5
$ws.Range("M2").Value = -1128.7894
$ws.Range("N2").Value = -1658.5
$ws.Range("H28").Value = 18997.25
$ws.Range("I28").Value = 15156.667
$ws.Range("K28").Value = 15156.667
$ws.Range("M28").Value = -14964.667
$ws.Range("H32").Value = 4848.6055
$ws.Range("I32").Value = 3840.0857
$ws.Range("J32").Value = 16614.666
$ws.Range("K32").Value = 3840.0857
$ws.Range("L32").Value = 16614.666
$ws.Range("M32").Value = -3553.0857
$ws.Range("N32").Value = -17188.666
$ws.Range("H74").Value = 2197.9285
$ws.Range("I74").Value = 1073
$ws.Range("J74").Value = 3322.8572
$ws.Range("K74").Value = 1073
$ws.Range("L74").Value = 3322.8572
$ws.Range("M74").Value = -199
$ws.Range("N74").Value = -5070.8572
$ws.Range("H77").Value = 2197.9285
$ws.Range("I77").Value = 1073
$ws.Range("J77").Value = 3322.8572
$ws.Range("K77").Value = 5365
$ws.Range("L77").Value = 16614.286
$ws.Range("M77").Value = -997
$ws.Range("N77").Value = -25350.286
$ws.Range("H99").Value = 18997.25
$ws.Range("I99").Value = 15156.667
$ws.Range("K99").Value = 15156.667
$ws.Range("M99").Value = -12161.667
$ws.Range("H110").Value = 2200.4614
$ws.Range("I110").Value = 1400.625
$ws.Range("K110").Value = 1400.625
$ws.Range("M110").Value = 644.375
$ws.Range("H116").Value = 1287.56
$ws.Range("I116").Value = 1241.7894
$ws.Range("J116").Value = 1432.5
$ws.Range("K116").Value = 1241.7894
$ws.Range("L116").Value = 1432.5
$ws.Range("M116").Value = 1052.2106
$ws.Range("N116").Value = -6020.5
$ws.Range("H132").Value = 2316.0833
$ws.Range("I132").Value = 1956.6
$ws.Range("J132").Value = 4113.5
$ws.Range("K132").Value = 5869.799999999999
$ws.Range("L132").Value = 12340.5
$ws.Range("M132").Value = -3339.799999999999
$ws.Range("N132").Value = -17400.5
$ws.Range("H133").Value = 97630.5
$ws.Range("J133").Value = 97630.5
$ws.Range("L133").Value = 97630.5
$ws.Range("N133").Value = -102690.5

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1287.56
$ws.Range("I3").Value = 1241.7894
$ws.Range("J3").Value = 1432.5
$ws.Range("K3").Value = 1241.7894
$ws.Range("L3").Value = 1432.5
$ws.Range("M3").Value = -1127.7894
$ws.Range("N3").Value = -1660.5
$ws.Range("H137").Value = 68721.75
$ws.Range("I137").Value = 7500
$ws.Range("J137").Value = 89129
$ws.Range("K137").Value = 7500
$ws.Range("L137").Value = 89129
$ws.Range("M137").Value = -2400
$ws.Range("N137").Value = -99329
$ws.Range("H138").Value = 77733.336
$ws.Range("I138").Value = 74950
$ws.Range("J138").Value = 100000
$ws.Range("K138").Value = 74950
$ws.Range("L138").Value = 100000
$ws.Range("M138").Value = -69810
$ws.Range("N138").Value = -110280

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1742.5555
$ws.Range("I58").Value = 1292.3
$ws.Range("J58").Value = 3029
$ws.Range("K58").Value = 1292.3
$ws.Range("L58").Value = 3029
$ws.Range("M58").Value = -1089.3
$ws.Range("N58").Value = -3435
$ws.Range("H132").Value = 3690.611
$ws.Range("I132").Value = 3741.8462
$ws.Range("J132").Value = 3557.4
$ws.Range("K132").Value = 11225.5386
$ws.Range("L132").Value = 10672.2
$ws.Range("M132").Value = -8695.5386
$ws.Range("N132").Value = -15732.2
$ws.Range("H136").Value = 1742.5555
$ws.Range("I136").Value = 1292.3
$ws.Range("J136").Value = 3029
$ws.Range("K136").Value = 3876.9
$ws.Range("L136").Value = 9087
$ws.Range("M136").Value = -1326.9
$ws.Range("N136").Value = -14187

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 71517.66
$ws.Range("J129").Value = 2947.4348
$ws.Range("L129").Value = 8842.304400000001
$ws.Range("N129").Value = -18842.3044
$ws.Range("H131").Value = 1832.9667
$ws.Range("I131").Value = 1644.875
$ws.Range("J131").Value = 1901.3636
$ws.Range("K131").Value = 4934.625
$ws.Range("L131").Value = 5704.0908
$ws.Range("M131").Value = 105.375
$ws.Range("N131").Value = -15784.0908
$ws.Range("H140").Value = 2270
$ws.Range("I140").Value = 2078.5715
$ws.Range("J140").Value = 4950
$ws.Range("K140").Value = 6235.7145
$ws.Range("L140").Value = 14850
$ws.Range("M140").Value = -1055.7145
$ws.Range("N140").Value = -25210
$ws.Range("H141").Value = 2971.6875
$ws.Range("I141").Value = 2636.4666
$ws.Range("K141").Value = 7909.399800000001
$ws.Range("M141").Value = -2729.399800000001

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 38433.168
$ws.Range("J46").Value = 46666
$ws.Range("L46").Value = 46666
$ws.Range("N46").Value = -46978
$ws.Range("H57").Value = 17073.2
$ws.Range("J57").Value = 26999.6
$ws.Range("L57").Value = 26999.6
$ws.Range("N57").Value = -28639.6
$ws.Range("H80").Value = 3388.8125
$ws.Range("I80").Value = 2730.9092
$ws.Range("K80").Value = 2730.9092
$ws.Range("M80").Value = -1732.9092
$ws.Range("H83").Value = 3388.8125
$ws.Range("I83").Value = 2730.9092
$ws.Range("K83").Value = 13654.546
$ws.Range("M83").Value = -8662.546
$ws.Range("H126").Value = 2558.2144
$ws.Range("I126").Value = 2632.6
$ws.Range("K126").Value = 7897.799999999999
$ws.Range("M126").Value = -5427.799999999999

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3424.3438
$ws.Range("I22").Value = 3469.3333
$ws.Range("J22").Value = 2749.5
$ws.Range("K22").Value = 3469.3333
$ws.Range("L22").Value = 2749.5
$ws.Range("M22").Value = -3174.3333
$ws.Range("N22").Value = -3339.5
$ws.Range("H27").Value = 3424.3438
$ws.Range("I27").Value = 3469.3333
$ws.Range("J27").Value = 2749.5
$ws.Range("K27").Value = 3469.3333
$ws.Range("L27").Value = 2749.5
$ws.Range("M27").Value = -3362.3333
$ws.Range("N27").Value = -2963.5
$ws.Range("H46").Value = 3817.35
$ws.Range("I46").Value = 4031.9092
$ws.Range("J46").Value = 3555.111
$ws.Range("K46").Value = 4031.9092
$ws.Range("L46").Value = 3555.111
$ws.Range("M46").Value = -3843.9092
$ws.Range("N46").Value = -3931.111
$ws.Range("H61").Value = 13703.419
$ws.Range("I61").Value = 2523.842
$ws.Range("K61").Value = 2523.842
$ws.Range("M61").Value = -2321.842
$ws.Range("H113").Value = 13703.419
$ws.Range("I113").Value = 2523.842
$ws.Range("K113").Value = 2523.842
$ws.Range("M113").Value = -353.8420000000001
$ws.Range("H122").Value = 70773.3
$ws.Range("I122").Value = 94679.27
$ws.Range("K122").Value = 284037.81
$ws.Range("M122").Value = -281587.81
$ws.Range("H132").Value = 4505.032
$ws.Range("I132").Value = 3679.32
$ws.Range("J132").Value = 7945.5
$ws.Range("K132").Value = 11037.96
$ws.Range("L132").Value = 23836.5
$ws.Range("M132").Value = -8507.960000000001
$ws.Range("N132").Value = -28896.5

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 52500
$ws.Range("J62").Value = 100000
$ws.Range("L62").Value = 100000
$ws.Range("N62").Value = -101248
$ws.Range("H65").Value = 52500
$ws.Range("J65").Value = 100000
$ws.Range("L65").Value = 500000
$ws.Range("N65").Value = -506240
$ws.Range("H108").Value = 98748.5
$ws.Range("J108").Value = 98748.5
$ws.Range("L108").Value = 98748.5
$ws.Range("N108").Value = -106428.5
$ws.Range("H113").Value = 455.16666
$ws.Range("I113").Value = 442.25
$ws.Range("J113").Value = 481
$ws.Range("K113").Value = 1326.75
$ws.Range("L113").Value = 1443
$ws.Range("M113").Value = 843.25
$ws.Range("N113").Value = -5783
$ws.Range("H132").Value = 2805.9546
$ws.Range("I132").Value = 1787.871
$ws.Range("K132").Value = 5363.613
$ws.Range("M132").Value = -2833.613

